$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 129 (shifts existing rows 129..171 down to 130..172)
$ws.Rows("129:129").Insert(-4121)  # xlShiftDown

# Copy the constant columns from what is now row 130 (the old row 129 data)
# into the new row 129, then set the row-specific new values.
$ws.Range("A129").Value = 9
$ws.Range("B129").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C129").Value = "Metropolitana"
$ws.Range("D129").Value = 45119
$ws.Range("E129").Value = 13
$ws.Range("F129").Value = 100112022
$ws.Range("G129").Value = "Arveja Verde"
$ws.Range("H129").Value = "Perfection"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 48
$ws.Range("K129").Value = 31000
$ws.Range("L129").Value = 33000
$ws.Range("M129").Value = 31917
$ws.Range("N129").Value = '$/malla 25 kilos'
$ws.Range("O129").Value = "Provincia de Limarí"
$ws.Range("P129").Value = 1277
$ws.Range("Q129").Value = 25
$ws.Range("R129").Value = "Hortaliza"

$ws.Range("D129").NumberFormat = $ws.Range("D130").NumberFormat
